$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from AC1 onto the new header cells, then set their values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 71
    $ws.Cells.Item($row, 31).Value = 91
    $ws.Cells.Item($row, 32).Value = 0
}
